# Update the dSF column (column F) values on the active worksheet.
# These reflect repulled data for the hendricks_kyle save file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = 2
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = -7
$ws.Range("F13").Value = -3
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = -3
$ws.Range("F17").Value = -2
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = -2
$ws.Range("F20").Value = 5
$ws.Range("F21").Value = 1
$ws.Range("F22").Value = 4
$ws.Range("F23").Value = -4
$ws.Range("F25").Value = -3
